{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block, along\n// with the blank paragraph that separated it from the bibliography entry.\n//\n// Target text markers (unique within the document):\n//   - the blank paragraph right after \"23.ed. S\u00e3o Paulo: Cortez, 2009.\"\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//     pages. Original theme under Creative Commons Attribution\"\n// The blank paragraph that follows the copyright line (right before the\n// page-break paragraph) is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the bibliography closing paragraph that marks where the footer\n// block begins (the paragraph right after it is the blank separator we\n// also need to remove).\nlet cortezIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"23.ed. S\u00e3o Paulo: Cortez, 2009.\") !== -1) {\n    cortezIndex = i;\n  } else if (text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  } else if (text.indexOf(\"Powered by Jekyll and Github pages\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (cortezIndex === -1 || jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\n    \"Could not locate one of the expected paragraphs (cortez=\" +\n      cortezIndex +\n      \", jupiter=\" +\n      jupiterIndex +\n      \", copyright=\" +\n      copyrightIndex +\n      \")\"\n  );\n}\n\n// The blank separator paragraph sits right between the bibliography\n// closing paragraph and the \"Ver no Jupiter\" paragraph.\nconst blankIndex = cortezIndex + 1;\nif (blankIndex !== jupiterIndex - 1) {\n  throw new Error(\n    \"Unexpected document layout: blank paragraph not immediately before the Jupiter paragraph\"\n  );\n}\n\n// Delete from bottom to top so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[jupiterIndex].delete();\nitems[blankIndex].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block, along\n# with the blank paragraph that separated it from the bibliography entry.\n#\n# Target text markers (unique within the document):\n#   - the blank paragraph right after \"23.ed. S\u00e3o Paulo: Cortez, 2009.\"\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#     pages. Original theme under Creative Commons Attribution\"\n# The blank paragraph that follows the copyright line (right before the\n# page-break paragraph) is left untouched.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$cortezIdx = -1\n$jupiterIdx = -1\n$copyrightIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t -like \"*Cortez, 2009.*\") { $cortezIdx = $i }\n  elseif ($t -like \"*Ver no Jupiter*\") { $jupiterIdx = $i }\n  elseif ($t -like \"*Powered by Jekyll*\") { $copyrightIdx = $i }\n}\n\nif ($cortezIdx -eq -1 -or $jupiterIdx -eq -1 -or $copyrightIdx -eq -1) {\n  throw \"Could not locate one of the expected paragraphs (cortez=$cortezIdx, jupiter=$jupiterIdx, copyright=$copyrightIdx)\"\n}\n\n# The blank separator paragraph sits right between the bibliography\n# closing paragraph and the \"Ver no Jupiter\" paragraph.\n$blankIdx = $cortezIdx + 1\nif ($blankIdx -ne ($jupiterIdx - 1)) {\n  throw \"Unexpected document layout: blank paragraph not immediately before the Jupiter paragraph\"\n}\n\n# Delete from bottom to top so earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIdx).Range.Delete()\n$d.Paragraphs.Item($jupiterIdx).Range.Delete()\n$d.Paragraphs.Item($blankIdx).Range.Delete()\n"}
